$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.907.79'
$ws.Range('E2').Value = '  +1.16%  '
$ws.Range('D3').Value = '2.488.72'
$ws.Range('E3').Value = '  +0.42%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '588.40'
$ws.Range('E5').Value = '  +0.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '174.28'
$ws.Range('E6').Value = '  +0.82%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').Value = '  +3.47%  '
$ws.Range('E10').Value = '  -1.59%  '
$ws.Range('E11').Value = '  +0.49%  '
$ws.Range('E12').Value = '  -0.09%  '
$ws.Range('E13').Value = '  +0.23%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.36'
$ws.Range('E14').Value = '  -0.51%  '
$ws.Range('D15').Value = '67.752.40'
$ws.Range('E15').Value = '  +0.99%  '
$ws.Range('E16').Value = '  -0.40%  '
$ws.Range('D17').Value = '2.461.93'
$ws.Range('E17').Value = '  +0.41%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.84'
$ws.Range('E18').Value = '  -1.03%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.36'
$ws.Range('E19').Value = '  -2.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '346.84'
$ws.Range('E20').Value = '  -0.90%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.12'
$ws.Range('E21').Value = '  +2.51%  '
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.81'
$ws.Range('E23').Value = '  +2.74%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.16'
$ws.Range('E24').Value = '  -1.29%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.71'
$ws.Range('E25').Value = '  -5.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.85'
$ws.Range('E26').Value = '  -2.84%  '
$ws.Range('E27').Value = '  +0.37%  '
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('D29').Value = '0.0₃0889'
$ws.Range('E29').Value = '  -2.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '497.76'
$ws.Range('E30').Value = '  -1.27%  '
$ws.Range('E31').Value = '  +0.52%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.25'
$ws.Range('E32').Value = '  +0.52%  '
$ws.Range('E33').Value = '  -0.22%  '
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '164.26'
$ws.Range('E35').Value = '  +1.18%  '
$ws.Range('E36').Value = '  +2.00%  '
$ws.Range('E37').Value = '  -0.37%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.21'
$ws.Range('E38').Value = '  +0.31%  '
$ws.Range('E40').Value = '  -1.79%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.73'
$ws.Range('E41').Value = '  +2.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.325'
$ws.Range('E42').Value = '  -0.91%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.77'
$ws.Range('E43').Value = '  -0.83%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.37'
$ws.Range('E44').Value = '  -0.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '148.40'
$ws.Range('E45').Value = '  +3.71%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.53'
$ws.Range('E46').Value = '  +1.49%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.512'
$ws.Range('E47').Value = '  -0.49%  '
$ws.Range('E48').Value = '  -4.58%  '
$ws.Range('E49').Value = '  -0.70%  '
$ws.Range('E50').Value = '  -1.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.575'
$ws.Range('E51').Value = '  -1.42%  '
